$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("G5").Value = 1.53
$ws.Range("H5").Value = 3.9
$ws.Range("I5").Value = 6.25
$ws.Range("J5").Value = 2.2
$ws.Range("L5").Value = 7
$ws.Range("U5").Value = 2.5
$ws.Range("V5").Value = 1.5
$ws.Range("X5").Value = 6
$ws.Range("Z5").Value = 10
$ws.Range("AD5").Value = 8
$ws.Range("AE5").Value = 26
$ws.Range("AF5").Value = 101
$ws.Range("AJ5").Value = 21
$ws.Range("AK5").Value = 81
$ws.Range("AM5").Value = 67
$ws.Range("AN5").Value = 3.25
$ws.Range("AO5").Value = 8
$ws.Range("AQ5").Value = 26
$ws.Range("AS5").Value = 251
$ws.Range("AX5").Value = 41

# Row 8 updates
$ws.Range("G8").Value = 3.45
$ws.Range("I8").Value = 1.95
$ws.Range("J8").Value = 3.9
$ws.Range("K8").Value = 2.18
$ws.Range("L8").Value = 2.5
$ws.Range("M8").Value = 8.2
$ws.Range("N8").Value = 1.06
$ws.Range("O8").Value = 1.29
$ws.Range("P8").Value = 3
$ws.Range("Q8").Value = 1.91
$ws.Range("R8").Value = 1.8
$ws.Range("V8").Value = 1.87
$ws.Range("W8").Value = 10.25
$ws.Range("X8").Value = 18
$ws.Range("Y8").Value = 12
$ws.Range("Z8").Value = 45
$ws.Range("AA8").Value = 32
$ws.Range("AB8").Value = 40
$ws.Range("AH8").Value = 7.2
$ws.Range("AI8").Value = 9.25
$ws.Range("AJ8").Value = 8.5
$ws.Range("AK8").Value = 17
$ws.Range("AL8").Value = 15.5
$ws.Range("AM8").Value = 28
$ws.Range("AN8").Value = 5.3
$ws.Range("AO8").Value = 18.5
$ws.Range("AP8").Value = 25
$ws.Range("AQ8").Value = 90
$ws.Range("AT8").Value = 2.65
$ws.Range("AU8").Value = 7.1
$ws.Range("AV8").Value = 60
$ws.Range("AW8").Value = 3.85
$ws.Range("AX8").Value = 9.5
$ws.Range("AY8").Value = 18
$ws.Range("AZ8").Value = 35
$ws.Range("BA8").Value = 65
$ws.Range("BB8").Value = 200

# Row 12 updates
$ws.Range("G12").Value = 2.05
$ws.Range("I12").Value = 3.6
$ws.Range("M12").Value = 1.07
$ws.Range("N12").Value = 8.5
$ws.Range("Q12").Value = 2.08
$ws.Range("R12").Value = 1.73
$ws.Range("U12").Value = 1.83
$ws.Range("V12").Value = 1.83
$ws.Range("Z12").Value = 19
$ws.Range("AC12").Value = 8.5
$ws.Range("AH12").Value = 10
$ws.Range("AL12").Value = 29
$ws.Range("AP12").Value = 23
